# fix: typos on menu
#
# Bokmål -> Nynorsk-style spelling fixes in the sandwich descriptions:
#   "håndstekt" -> "handsteikt", "stekt" -> "steikt", "løk" -> "lauk"
# plus the selected/active cell on the sheet moved from A5 to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "Seilerpatent" description (column C)
$ws.Range("C4").Value = "Handsteikt karbonade m/ speilegg. Havrebrød, salat, tomat og agurk. 1,2,3,8,10"

# Row 3: "Karbonade m/ løk" description (column C)
$ws.Range("C3").Value = "Havrebrød, handsteikt karbonade med steikt lauk, salat, agurk, tomat, sylteagurk og svisker. 1,2,8,10"

# Active selection moved to C5
$ws.Range("C5").Select()
